# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (H) and "Correspond Handback
# DateTime" (K) values for the first data row (4325be50-...md) on the
# zh-cn and de-de language sheets, reflecting a fresh handback report run.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("H2").Value = "2016-08-17 10:44:44"
$ws_zhcn.Range("K2").Value = "2016-08-17 10:45:20"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("H2").Value = "2016-08-17 10:44:51"
$ws_dede.Range("K2").Value = "2016-08-17 10:45:27"
